# Big-Picture Style Notification implementieren
# Add a new time-tracking entry (row 8) to the Arbeitszeit worksheet:
#   01.11.2021 | 2h | "Praesentation, Aussehen, Big-Picture"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Reuse the formatting of the row above (border + number format) by
# copying it down, then overwrite the copied values with the real ones.
$ws.Range("A7:B7").Copy()
$ws.Range("A8:B8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A8").Value = Get-Date -Year 2021 -Month 11 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Präsentation, Aussehen, Big-Picture"

$ws.Range("A8:B8").Select()
